$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# New porcelanato sale recorded for FUENTES PAREDES MARIA FERNANDA (row 12)
$ws1.Range("M12").Value = 333.59

# New inodoros / lavabos sales recorded for MANCHENO PINO HERVIN SANTIAGO (row 20)
$ws1.Range("H20").Value = 71.09999999999999
$ws1.Range("I20").Value = 26.1

# Updated "x de 34" progress counters on the totals row (row 36)
$ws1.Range("H36").Value = "1 de 34"
$ws1.Range("I36").Value = "1 de 34"
$ws1.Range("M36").Value = "5 de 34"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# "octubre" column totals updated to reflect the new sales above
$ws2.Range("F12").Value = 523.67
$ws2.Range("F20").Value = 1515.39
$ws2.Range("F36").Value = 9716.360000000001

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# INODOROS row
$ws3.Range("D6").Value = 71.09999999999999
$ws3.Range("E6").Value = 743.0234308088729
$ws3.Range("F6").Value = 0.0873331945861803

# LAVABOS row
$ws3.Range("D7").Value = 26.1
$ws3.Range("E7").Value = 195.577754071894
$ws3.Range("F7").Value = 0.1177384718158743

# PORCELANATO row
$ws3.Range("D12").Value = 8405.549999999999
$ws3.Range("E12").Value = 13295.72
$ws3.Range("F12").Value = 0.3873298659479376

# TOTAL row
$ws3.Range("D14").Value = 9716.359999999999
$ws3.Range("E14").Value = 26869.20723718182
$ws3.Range("F14").Value = 0.2655790447913374
